$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("individuals")
$ws.Select()

$ws.Range("BN2").Value = "a"
$ws.Range("BV2").Value = "b"

$ws.Range("BK1").Interior.Pattern = -4142  # xlNone
$ws.Range("BL1").Interior.Pattern = -4142  # xlNone
$ws.Range("BM1").Interior.ThemeColor = 5
$ws.Range("BN1").Interior.Color = 15773696
$ws.Range("BO1").Interior.Color = 65535
$ws.Range("BP1").Interior.Color = 10498160
$ws.Range("BQ1").Interior.ThemeColor = 7
$ws.Range("BR1").Interior.ThemeColor = 9
$ws.Range("BS1").Interior.ThemeColor = 7
$ws.Range("BT1").Interior.Color = 65535
$ws.Range("BU1").Interior.Color = 10498160
$ws.Range("BV1").Interior.Color = 15773696

$ws.Range("BN2").Select()
